$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing data runs from row 2 through row 328 (date serials in column A,
# counts in B/C/D). Extend the series with 15 more days, from serial 44403
# (2021-07-26) through serial 44417 (2021-08-09) - "aggiornamento fino a 9
# agosto 2021" - replicating the same all-zero values and the date-cell
# formatting used by the existing rows.

$startRow = 329
$startSerial = 44403
$endSerial = 44417
$templateRow = 328

for ($serial = $startSerial; $serial -le $endSerial; $serial++) {
    $row = $startRow + ($serial - $startSerial)

    $ws.Cells.Item($row, 1).Value = $serial
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0

    # Clone the date cell's style (border/font/alignment/number format) from
    # the last existing row rather than just the number format, so the new
    # cell lands on the same shared style index as the rest of column A.
    $ws.Range("A" + $templateRow).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
}
